$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "HFE / Huainan, China" row (row 261). Deleting the entire row
# shifts all subsequent rows up by one, which also accounts for the
# dimension shrinking from A1:H334 to A1:H333.
$ws.Rows(261).Delete()
